# Increment columns E (taxa) and F (outra taxa) by a factor of 100
# for the ranking summary table, rows 2 through 16.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 16; $row++) {
    foreach ($col in @("E", "F")) {
        $cell = $ws.Range("$col$row")
        $cell.Value = $cell.Value2 * 100
    }
}
